$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.885583758354187
$ws.Range("B1").Value = 1.381412029266357
$ws.Range("C1").Value = 4.767478942871094
$ws.Range("D1").Value = 3.281548261642456
$ws.Range("E1").Value = 0.4603109359741211
